$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @{
    54 = @("2024-05-04 21:08:45", 12, 8, 3, 2, 3, 0, 0, 0.001, 0.05, 0.003, 100, 500, 10, 9, 3, 400, 5, 1, 90, 0.6666666666666666, "Data/bombay1.xlsx", -41200, "No es Simulación")
    55 = @("2024-05-04 23:16:28", 14, 9, 0, 2, 2, 5, 0, 0.001, 0.05, 0.003, 100, 500, 10, 9, 3, 100, 5, 1, 90, 0.6428571428571429, "Data/bombay1.xlsx", -23900, "No es Simulación")
    56 = @("2024-05-06 21:50:51", 9, 5, 1, 2, 0, 2, 0, 0.001, 0.05, 0.003, 100, 500, 10, 9, 3, 100, 5, 1, 90, 0.5555555555555556, "Data/bombay1.xlsx", -45300, "No es Simulación")
    57 = @("2024-05-09 23:17:23", 10, 6, 1, 0, 1, 4, 0, 0.001, 0.05, 0.003, 100, 500, 10, 9, 3, 400, 5, 1, 90, 0.6, "Data/bombay1.xlsx", -73200, "No es Simulación")
    58 = @("2024-05-14 00:10:08", 7, 6, 0, 1, 3, 2, 0, 0.001, 0.05, 0.003, 100, 500, 10, 9, 3, 400, 5, 1, 90, 0.8571428571428571, "Data/bombay1.xlsx", 47200, "No es Simulación")
    59 = @("2024-05-26 01:14:41", 9, 7, 0, 1, 3, 3, 0, 0.001, 0.05, 0.003, 100, 500, 10, 9, 3, 1000, 5, 1, 90, 0.7777777777777778, "Data/bombay1.xlsx", 54000, "No es Simulación")
}

foreach ($r in ($newRows.Keys | Sort-Object)) {
    $rowValues = $newRows[$r]
    for ($i = 0; $i -lt $rowValues.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $rowValues[$i]
    }
}

Write-Host "Added rows 54-59"